# repo clean-up, jochen's slides
#
# Applies the content-level edits captured by the commit:
#   1. Fix the cached "datetimeFigureOut" footer field on the two slide
#      layouts that still show the old 2-digit-year rendering
#      ("20.05.13" -> "20.05.2013").
#   2. Add the missing "Matze!" body text on slide 4 (the other slides
#      already carry it).
#
$p = $ppt.ActivePresentation

# --- 1. Date placeholder fix-up on the slide layouts -----------------
$master = $p.SlideMaster

# slideLayout1.xml ("Titelfolie" / title layout) - shape "Rectangle 29"
$dateShape1 = $master.CustomLayouts.Item(1).Shapes.Item(6)
if ($dateShape1.TextFrame.TextRange.Text -eq "20.05.13") {
    $dateShape1.TextFrame.TextRange.Text = "20.05.2013"
}

# slideLayout10.xml - shape "Datumsplatzhalter 3"
$dateShape2 = $master.CustomLayouts.Item(10).Shapes.Item(3)
if ($dateShape2.TextFrame.TextRange.Text -eq "20.05.13") {
    $dateShape2.TextFrame.TextRange.Text = "20.05.2013"
}

# --- 2. Add "Matze!" to slide 4's empty content placeholder -----------
$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Text = "Matze!"
